# Fix the "Further reading" slide's readings line:
#   "J&M2, Ch 21.1-21.3 , 21.5-21.6"
# becomes
#   "J&M2, Ch 21.1-21.3, 21.5-21.6"
# and the fragment that used to be one run
#   " 21.1-21.3 , 21.5-21.6"
# is split into three separate runs (same rPr on each):
#   " " / "21.1-21.3, " / "21.5-21.6"

$p = $ppt.ActivePresentation

# Locate the slide/shape holding the "Ch 21.1-21.3" reading line instead of
# hard-coding indices, so the script is resilient to minor structural drift.
$targetSlide = $null
$targetShape = $null
$targetPara = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text.Contains("21.1-21.3")) {
                $targetSlide = $slide
                $targetShape = $shape
                $targetPara = $para
                break
            }
        }
        if ($targetPara -ne $null) { break }
    }
    if ($targetPara -ne $null) { break }
}

if ($targetPara -eq $null) {
    throw "could not locate the '21.1-21.3' reading paragraph"
}

$paraText = $targetPara.Text

# Old fragment (note the stray space before the comma) and its replacement.
$oldFragment = " 21.1-21.3 , 21.5-21.6"
$newFull     = " 21.1-21.3, 21.5-21.6"
$part1 = " "
$part2 = "21.1-21.3, "
$part3 = "21.5-21.6"

$fragStart0 = $paraText.IndexOf($oldFragment)
if ($fragStart0 -lt 0) {
    throw "could not locate the old reading fragment text"
}

# TextRange character positions are 1-based.
$fragStart = $fragStart0 + 1
$fragLen = $oldFragment.Length

# Step 1: replace the whole fragment in one go (keeps it as a single run for now).
$whole = $targetPara.Characters($fragStart, $fragLen)
$whole.Text = $newFull

# Step 2: re-carve the same span into three separate runs by assigning text to
# disjoint sub-ranges individually -- each assignment keeps its own run instead
# of being re-merged into neighbouring runs of identical formatting.
$r1 = $targetPara.Characters($fragStart, $part1.Length)
$r1.Text = $part1

$r2 = $targetPara.Characters($fragStart + $part1.Length, $part2.Length)
$r2.Text = $part2

$r3 = $targetPara.Characters($fragStart + $part1.Length + $part2.Length, $part3.Length)
$r3.Text = $part3

Write-Host "Updated paragraph text: [$($targetShape.Name) on slide $($targetSlide.SlideIndex)] -> $($targetPara.Text)"
